$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E22) listed periods in descending order
# (2507..2501). The update re-sorts this block into ascending chronological
# order (2501..2507) ahead of adding "parte 1" of the new account statements.
$ws.Range("E16").Value = "2501"
$ws.Range("E17").Value = "2502"
$ws.Range("E18").Value = "2503"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2505"
$ws.Range("E21").Value = "2506"
$ws.Range("E22").Value = "2507"

# The "Valor Mora" figures travel with their period row, so the values that
# belonged to period 2501 and 2507 swap places along with the labels above.
$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 48533
